# Rename the "_old"/"_new" column header suffixes to the respective
# input-file-version suffixes ("_FV2410" / "_FV2504"), then turn the
# used range into a native Excel Table ("Table1") and freeze the header
# row, matching the target workbook structure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J hold the "old" (FV2410) snapshot, columns L-U hold the
# "new" (FV2504) snapshot; column K is the untouched "diff" column.
$leftCols  = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$rightCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($leftCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
    $ws.Range($rightCols[$i] + "1").Value = $baseNames[$i] + "_FV2504"
}

# Turn the data range into a table so the header row gets AutoFilter +
# structured-table metadata (xl/tables/table1.xml).
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U77"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# Freeze the header row (split after row 1, keep column A as the
# left-most visible column of the scrolling pane).
$ws.Activate()
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
